$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "C"
$ws.Range("B3").Value = "D"
$ws.Range("B1").Select()
